# Auto-generated edit script applying the Kraken_Profits.xlsx diff
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns across 8 sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 35000
$ws.Range("I18").Value = 35000
$ws.Range("K18").Value = 35000
$ws.Range("M18").Value = -34716
$ws.Range("H80").Value = 8616.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 8616.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 25849.5
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -27845.5
$ws.Range("H83").Value = 8616.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 8616.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 77548.5
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -87532.5
$ws.Range("H100").Value = 6114.7144
$ws.Range("I100").Value = 5726
$ws.Range("K100").Value = 5726
$ws.Range("M100").Value = -5185
$ws.Range("H132").Value = 3496.56
$ws.Range("I132").Value = 2610
$ws.Range("J132").Value = 9998
$ws.Range("K132").Value = 7830
$ws.Range("L132").Value = 29994
$ws.Range("M132").Value = -5300
$ws.Range("N132").Value = -35054

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1160
$ws.Range("I45").Value = 1160
$ws.Range("K45").Value = 1160
$ws.Range("M45").Value = -783
$ws.Range("H101").Value = 12998
$ws.Range("J101").Value = 12998
$ws.Range("L101").Value = 12998
$ws.Range("N101").Value = -19488
$ws.Range("H110").Value = 1166.6666
$ws.Range("I110").Value = 1166.6666
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1166.6666
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 878.3334
$ws.Range("N110").ClearContents()
$ws.Range("H111").Value = 99995
$ws.Range("J111").Value = 99995
$ws.Range("L111").Value = 99995
$ws.Range("N111").Value = -108175
$ws.Range("H122").Value = 2497.5
$ws.Range("I122").Value = 2497.5
$ws.Range("K122").Value = 7492.5
$ws.Range("M122").Value = -5042.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1250
$ws.Range("I20").Value = 2000
$ws.Range("K20").Value = 2000
$ws.Range("M20").Value = -1753
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H80").Value = 520.4
$ws.Range("J80").Value = 803
$ws.Range("L80").Value = 803
$ws.Range("N80").Value = -2799
$ws.Range("H83").Value = 520.4
$ws.Range("J83").Value = 803
$ws.Range("L83").Value = 4015
$ws.Range("N83").Value = -13999
$ws.Range("H94").Value = 2351.3809
$ws.Range("J94").Value = 3126.25
$ws.Range("L94").Value = 3126.25
$ws.Range("N94").Value = -4028.25
$ws.Range("H100").Value = 9749.75
$ws.Range("J100").Value = 9749.75
$ws.Range("L100").Value = 9749.75
$ws.Range("N100").Value = -11913.75
$ws.Range("H103").Value = 3651.6667
$ws.Range("J103").Value = 3651.6667
$ws.Range("L103").Value = 3651.6667
$ws.Range("N103").Value = -5995.6667
$ws.Range("H105").Value = 6500
$ws.Range("I105").Value = 8000
$ws.Range("K105").Value = 8000
$ws.Range("M105").Value = -6253
$ws.Range("H134").Value = 6571.2856
$ws.Range("I134").Value = 5000
$ws.Range("K134").Value = 15000
$ws.Range("M134").Value = -12465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1266.6666
$ws.Range("I16").Value = 1266.6666
$ws.Range("K16").Value = 1266.6666
$ws.Range("M16").Value = -979.6666
$ws.Range("H43").Value = 28666.334
$ws.Range("J43").Value = 28666.334
$ws.Range("L43").Value = 28666.334
$ws.Range("N43").Value = -29034.334
$ws.Range("H101").Value = 28666.334
$ws.Range("J101").Value = 28666.334
$ws.Range("L101").Value = 28666.334
$ws.Range("N101").Value = -35156.334
$ws.Range("H105").Value = 3373.25
$ws.Range("I105").Value = 3247
$ws.Range("J105").Value = 3499.5
$ws.Range("K105").Value = 3247
$ws.Range("L105").Value = 3499.5
$ws.Range("M105").Value = -1500
$ws.Range("N105").Value = -6993.5
$ws.Range("H110").Value = 99995
$ws.Range("J110").Value = 99995
$ws.Range("L110").Value = 99995
$ws.Range("N110").Value = -108175
$ws.Range("H113").Value = 1266.6666
$ws.Range("I113").Value = 1266.6666
$ws.Range("K113").Value = 1266.6666
$ws.Range("M113").Value = 903.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 130
$ws.Range("I17").Value = 106.666664
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 319.999992
$ws.Range("L17").Value = 600
$ws.Range("M17").Value = -150.999992
$ws.Range("N17").Value = -938
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H92").Value = 375
$ws.Range("I92").Value = 350
$ws.Range("J92").Value = 400
$ws.Range("K92").Value = 1050
$ws.Range("L92").Value = 1200
$ws.Range("M92").Value = 198
$ws.Range("N92").Value = -3696
$ws.Range("H123").Value = 600
$ws.Range("I123").Value = 600
$ws.Range("K123").Value = 1800
$ws.Range("M123").Value = 650
$ws.Range("H131").Value = 2017.4615
$ws.Range("J131").Value = 2244.2222
$ws.Range("L131").Value = 6732.6666
$ws.Range("N131").Value = -16812.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4500
$ws.Range("J80").Value = 4500
$ws.Range("L80").Value = 4500
$ws.Range("N80").Value = -6496
$ws.Range("H83").Value = 4500
$ws.Range("J83").Value = 4500
$ws.Range("L83").Value = 22500
$ws.Range("N83").Value = -32484
$ws.Range("H122").Value = 2014.7142
$ws.Range("I122").Value = 2161
$ws.Range("J122").Value = 1649
$ws.Range("K122").Value = 6483
$ws.Range("L122").Value = 4947
$ws.Range("M122").Value = -4033
$ws.Range("N122").Value = -9847
$ws.Range("H132").Value = 8999.5
$ws.Range("I132").Value = 7999.5
$ws.Range("J132").Value = 9999.5
$ws.Range("K132").Value = 23998.5
$ws.Range("L132").Value = 29998.5
$ws.Range("M132").Value = -21468.5
$ws.Range("N132").Value = -35058.5
$ws.Range("H135").Value = 100000
$ws.Range("J135").Value = 100000
$ws.Range("L135").Value = 100000
$ws.Range("N135").Value = -110140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2338.6924
$ws.Range("J55").Value = 2386.75
$ws.Range("L55").Value = 2386.75
$ws.Range("N55").Value = -2732.75
$ws.Range("H132").Value = 8666.333000000001
$ws.Range("I132").Value = 10000
$ws.Range("J132").Value = 7999.5
$ws.Range("K132").Value = 30000
$ws.Range("L132").Value = 23998.5
$ws.Range("M132").Value = -27470
$ws.Range("N132").Value = -29058.5
$ws.Range("H136").Value = 5475.75
$ws.Range("I136").Value = 4952
$ws.Range("J136").Value = 5999.5
$ws.Range("K136").Value = 14856
$ws.Range("L136").Value = 17998.5
$ws.Range("M136").Value = -12306
$ws.Range("N136").Value = -23098.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 41349.7
$ws.Range("I4").Value = 67584
$ws.Range("J4").Value = 1998.25
$ws.Range("K4").Value = 67584
$ws.Range("L4").Value = 1998.25
$ws.Range("M4").Value = -67471
$ws.Range("N4").Value = -2224.25
$ws.Range("H122").Value = 2758.8572
$ws.Range("I122").Value = 2758.8572
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8276.571599999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5826.571599999999
$ws.Range("N122").ClearContents()

Write-Output "Applied $(193) value updates and $(6) clears across 8 sheets"